$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New row 6: a new test entry ("JEAN CARLOS") ---
$ws.Range("A6").Value = 5
$ws.Range("B6").Value = "JEAN CARLOS"
$ws.Range("C6").Value = "20/12/2022"
$ws.Range("D6").Value = "00:00"
$ws.Range("E6").Value = 50
$ws.Range("F6").Value = 50
$ws.Range("G6").Value = 50
$ws.Range("H6").Value = 50
$ws.Range("J6").Value = "TESTE`nTESTE`nTESTE`nTESTE`nTESTE`nTESE`nTESTE`nTESTE`nTESTE`nTEST`nTEST`nTES`nTE`nTEST`nTESTE`nETSTE`nTETSS"
$ws.Range("K6").Value = "Pendente"

# Let the row height auto-size to the wrapped/multi-line content, then let Excel
# settle back to the sheet's default (no explicit custom height is persisted).
$ws.Rows(6).EntireRow.AutoFit() | Out-Null

# --- C5 was retyped as plain text instead of a date ---
$ws.Range("C5").ClearContents() | Out-Null
$ws.Range("C5").NumberFormat = "@"
$ws.Range("C5").Value = "01/01/2023"
$ws.Range("C5").Style = "Normal"

# --- Column C was widened (no longer auto "best fit") ---
$ws.Columns("C").ColumnWidth = 17.59

# --- Cursor ended up on D9 ---
$ws.Range("D9").Select() | Out-Null
